$wb = $excel.ActiveWorkbook

# Rename the first sheet from "TestData" to "TestData1"
$ws1 = $wb.Worksheets.Item("TestData")
$ws1.Name = "TestData1"

# Move the selection on that sheet from D5 to F4
$ws1.Activate()
$ws1.Range("F4").Select()
